$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the ASIN column for the week-start date.
$ws.Columns.Item(2).Insert()
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week start dates for each week row (weekly cadence starting 2025-01-05).
$weekDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekDates.Length; $i++) {
    $row = $i + 2

    # Week labels drop the leading zero (W01 -> W1, ... W09 -> W9).
    $ws.Cells.Item($row, 1).Value = "W" + ($i + 1)

    # Write the date as literal text (not an Excel date serial number):
    # force a text number format before assigning, then clear the
    # formatting again so the cell keeps its default style but still
    # holds the string value.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $weekDates[$i]
    $dateCell.ClearFormats()
}
